$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 31536
$ws.Range("B2").Value = "Pedro Miguel Novaes"
$ws.Range("C2").Value = "Vendas"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45101
$ws.Range("G2").Value = 3606.19

# Row 3
$ws.Range("A3").Value = 68965
$ws.Range("B3").Value = "Gabriel Sales"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45089
$ws.Range("G3").Value = 3127.05

# Row 4
$ws.Range("A4").Value = 48398
$ws.Range("B4").Value = "Agatha da Rosa"
$ws.Range("C4").Value = "Jurídico"
$ws.Range("D4").Value = "Problemas pessoais"
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 45090
$ws.Range("G4").Value = 12481.63

# Row 5
$ws.Range("A5").Value = 11878
$ws.Range("B5").Value = "Letícia Ribeiro"
$ws.Range("C5").Value = "P&D"
$ws.Range("D5").Value = "Outros"
$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 45103
$ws.Range("G5").Value = 4673.88

# Row 6
$ws.Range("A6").Value = 41024
$ws.Range("B6").Value = "Beatriz Jesus"
$ws.Range("C6").Value = "Vendas"
$ws.Range("D6").Value = "Viagem de negócios"
$ws.Range("E6").Value = 7
$ws.Range("F6").Value = 45092
$ws.Range("G6").Value = 6540.34

# Row 7
$ws.Range("A7").Value = 39221
$ws.Range("B7").Value = "João Pedro Fogaça"
$ws.Range("C7").Value = "Atendimento ao Cliente"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 45101
$ws.Range("G7").Value = 3729.67

# Row 8
$ws.Range("A8").Value = 37621
$ws.Range("B8").Value = "Esther da Conceição"
$ws.Range("C8").Value = "Jurídico"
$ws.Range("D8").Value = "Doença"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45104
$ws.Range("G8").Value = 7467.38

# Row 9
$ws.Range("A9").Value = 24703
$ws.Range("B9").Value = "Joana Duarte"
$ws.Range("C9").Value = "Financeiro"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 45100
$ws.Range("G9").Value = 6641.13

# Row 10
$ws.Range("A10").Value = 65618
$ws.Range("B10").Value = "Marina Gomes"
$ws.Range("C10").Value = "P&D"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45090
$ws.Range("G10").Value = 10926.48

# Row 11
$ws.Range("A11").Value = 23842
$ws.Range("B11").Value = "Elisa Azevedo"
$ws.Range("C11").Value = "Engenharia"
$ws.Range("D11").Value = "Outros"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45078
$ws.Range("G11").Value = 8684.24
